# Generate Report for Handback
#
# The two tracked source files have come back from translation "in sync
# with en-US", so this refreshes the localization-status report:
#  - the status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" (Overview sheet, and the Status
#    column on each language sheet)
#  - each language sheet's "Latest Target File" / "Latest Handback File"
#    columns (E/F) get populated with the same file names as the source /
#    handoff columns (A/C), as hyperlinks
#  - the "Latest Handback DateTime" column (G) gets the timestamp of the
#    handback

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# Blue used by the workbook's built-in "HyperLink" cell style (FF6495ED),
# expressed the way VBA/OOXML store Font.Color (0x00BBGGRR).
$hyperlinkColor = 15570276

# ---- Overview sheet: refresh the status text for both tracked files ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value2 = $statusNew
$overview.Range("C2").Value2 = $statusNew
$overview.Range("B3").Value2 = $statusNew
$overview.Range("C3").Value2 = $statusNew

# ---- Every language sheet's Status column shares the same text ----
$zhcnStatus = $wb.Worksheets.Item("zh-cn")
$zhcnStatus.Range("B2").Value2 = $statusNew
$zhcnStatus.Range("B3").Value2 = $statusNew

$dedeStatus = $wb.Worksheets.Item("de-de")
$dedeStatus.Range("B2").Value2 = $statusNew
$dedeStatus.Range("B3").Value2 = $statusNew

function Update-LangSheet($SheetName, $HandbackTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Row 2 -> 54dead1e-3150-4dbd-bb1b-e2d47481260c.md
    $row2Name = $ws.Range("A2").Value2
    $row2HandoffFile = $ws.Range("C2").Value2
    $row2HandoffUrl = $ws.Hyperlinks.Item(2).Address

    $ws.Range("E2").Value2 = $row2Name
    $ws.Hyperlinks.Add($ws.Range("E2"), $ws.Hyperlinks.Item(1).Address, [Type]::Missing, [Type]::Missing, $row2Name) | Out-Null
    $ws.Range("F2").Value2 = $row2HandoffFile
    $ws.Hyperlinks.Add($ws.Range("F2"), $row2HandoffUrl, [Type]::Missing, [Type]::Missing, $row2HandoffFile) | Out-Null
    $ws.Range("G2").Value2 = $HandbackTime

    # Row 3 -> 9ba2b223-8801-4341-9a8c-da5d7a438d89.md
    $row3Name = $ws.Range("A3").Value2
    $row3HandoffFile = $ws.Range("C3").Value2
    $row3HandoffUrl = $ws.Hyperlinks.Item(4).Address

    $ws.Range("E3").Value2 = $row3Name
    $ws.Hyperlinks.Add($ws.Range("E3"), $ws.Hyperlinks.Item(3).Address, [Type]::Missing, [Type]::Missing, $row3Name) | Out-Null
    $ws.Range("F3").Value2 = $row3HandoffFile
    $ws.Hyperlinks.Add($ws.Range("F3"), $row3HandoffUrl, [Type]::Missing, [Type]::Missing, $row3HandoffFile) | Out-Null
    $ws.Range("G3").Value2 = $HandbackTime

    # Match the hyperlink look (blue + underline) already used by A/C
    $newCells = $ws.Range("E2:F3")
    $newCells.Font.Underline = $true
    $newCells.Font.Color = $hyperlinkColor
}

Update-LangSheet "zh-cn" "2016-03-09 14:15:27"
Update-LangSheet "de-de" "2016-03-09 14:15:32"
